# Added fun INDEX implementation
# Mirrors the "Range with holes / VLOOKUP" demo block (rows 29-41) with a
# new "INDEX" demo block appended below it (rows 43-53).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Formatting -----------------------------------------------------
# Reuse the existing header style (row 29, "Range with holes") for the
# new "INDEX" header row, and the existing striped-fill style used by
# the static-fields rows (30:32) for the new static-fields rows (44:46).
$ws.Range("A29:D29").Copy() | Out-Null
$ws.Range("A43:D43").PasteSpecial(-4122) | Out-Null

$ws.Range("A30:D32").Copy() | Out-Null
$ws.Range("A44:D46").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Section header ---------------------------------------------------
$ws.Range("A43").Value = "INDEX"

# --- Static data table (B44:D46) --------------------------------------
$ws.Range("A44").Value = "static fields"
$ws.Range("B44").Value = 1
$ws.Range("C44").Value = 2
$ws.Range("D44").Value = 3

$ws.Range("A45").Value = "static fields"
$ws.Range("B45").Value = 10
$ws.Range("C45").Value = 20
$ws.Range("D45").Value = 30

$ws.Range("A46").Value = "static fields"
$ws.Range("B46").Value = 100
$ws.Range("C46").Value = 200
$ws.Range("D46").Value = 300

# --- FORMULATEXT column (A47:A53), entered once and filled down so it
# becomes a single shared formula, exactly like the VLOOKUP block above.
$ws.Range("A47:A53").Formula = '=FORMULATEXT(B47)'

# --- INDEX() formulas (B47:B53) ----------------------------------------
$ws.Range("B47").Formula = '=INDEX($B$44:$D$44,0)'
$ws.Range("B48").Formula = '=INDEX($B$44:$D$44,1)'
$ws.Range("B49").Formula = '=INDEX($B$44:$D$44,2)'
$ws.Range("B50").Formula = '=INDEX($B$44:$B$46,1)'
$ws.Range("B51").Formula = '=INDEX($B$44:$B$46,2)'
$ws.Range("B52").Formula = '=INDEX($B$44:$D$46,1)'
$ws.Range("B53").Formula = '=INDEX($B$44:$D$46,2)'

# --- View state: scroll / selection moved down onto the new block ------
$ws.Activate() | Out-Null
$excel.Goto($ws.Range("A27"), $true) | Out-Null
$ws.Range("B48").Select() | Out-Null
